{"js": "const replacements = [\n  [\"2024-03-22 Friday\", \"2024-03-23 Saturday\"],\n  [\"255\\u00d74=1020\", \"160\\u00d76=960\"],\n  [\"267\\u00d74=1068\", \"342\\u00d79=3078\"],\n  [\"217\\u00d76=1302\", \"415\\u00d78=3320\"],\n  [\"954\\u00d79=8586\", \"856\\u00d74=3424\"],\n  [\"529\\u00d76=3174\", \"731\\u00d74=2924\"],\n  [\"261\\u00d75=1305\", \"424\\u00d72=848\"],\n  [\"672\\u00d73=2016\", \"111\\u00d74=444\"],\n  [\"113\\u00d73=339\", \"863\\u00d75=4315\"],\n  [\"318\\u00d78=2544\", \"110\\u00d76=660\"],\n  [\"222\\u00d76=1332\", \"992\\u00d79=8928\"],\n  [\"345\\u00d74=1380\", \"245\\u00d72=490\"],\n  [\"101\\u00d75=505\", \"447\\u00d75=2235\"],\n  [\"198\\u00d75=990\", \"980\\u00d73=2940\"],\n  [\"846\\u00d77=5922\", \"750\\u00d77=5250\"],\n  [\"723\\u00d78=5784\", \"726\\u00d78=5808\"],\n  [\"580\\u00d76=3480\", \"349\\u00d79=3141\"],\n  [\"484\\u00d73=1452\", \"118\\u00d79=1062\"],\n  [\"394\\u00d78=3152\", \"579\\u00d72=1158\"],\n  [\"177\\u00d79=1593\", \"481\\u00d79=4329\"],\n  [\"609\\u00d75=3045\", \"875\\u00d78=7000\"],\n  [\"202\\u00d74=808\", \"983\\u00d72=1966\"],\n  [\"206\\u00d75=1030\", \"269\\u00d79=2421\"],\n  [\"212\\u00d76=1272\", \"423\\u00d72=846\"],\n  [\"231\\u00d77=1617\", \"509\\u00d74=2036\"],\n  [\"204\\u00d72=408\", \"437\\u00d74=1748\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-22 Friday\", \"2024-03-23 Saturday\"),\n    @(\"255\u00d74=1020\", \"160\u00d76=960\"),\n    @(\"267\u00d74=1068\", \"342\u00d79=3078\"),\n    @(\"217\u00d76=1302\", \"415\u00d78=3320\"),\n    @(\"954\u00d79=8586\", \"856\u00d74=3424\"),\n    @(\"529\u00d76=3174\", \"731\u00d74=2924\"),\n    @(\"261\u00d75=1305\", \"424\u00d72=848\"),\n    @(\"672\u00d73=2016\", \"111\u00d74=444\"),\n    @(\"113\u00d73=339\", \"863\u00d75=4315\"),\n    @(\"318\u00d78=2544\", \"110\u00d76=660\"),\n    @(\"222\u00d76=1332\", \"992\u00d79=8928\"),\n    @(\"345\u00d74=1380\", \"245\u00d72=490\"),\n    @(\"101\u00d75=505\", \"447\u00d75=2235\"),\n    @(\"198\u00d75=990\", \"980\u00d73=2940\"),\n    @(\"846\u00d77=5922\", \"750\u00d77=5250\"),\n    @(\"723\u00d78=5784\", \"726\u00d78=5808\"),\n    @(\"580\u00d76=3480\", \"349\u00d79=3141\"),\n    @(\"484\u00d73=1452\", \"118\u00d79=1062\"),\n    @(\"394\u00d78=3152\", \"579\u00d72=1158\"),\n    @(\"177\u00d79=1593\", \"481\u00d79=4329\"),\n    @(\"609\u00d75=3045\", \"875\u00d78=7000\"),\n    @(\"202\u00d74=808\", \"983\u00d72=1966\"),\n    @(\"206\u00d75=1030\", \"269\u00d79=2421\"),\n    @(\"212\u00d76=1272\", \"423\u00d72=846\"),\n    @(\"231\u00d77=1617\", \"509\u00d74=2036\"),\n    @(\"204\u00d72=408\", \"437\u00d74=1748\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
